$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing row 13 values (C13: 13.0 -> 14.0, D13: 57.0 -> 60.0)
$ws.Range("C13").Value = 14.0
$ws.Range("D13").Value = 60.0

# Add new row 17 - "pay bills" section entry
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "feemicon pill"
$ws.Range("C17").Value = 69.69
$ws.Range("D17").Value = 10000

$wb.Save()
